# "Generate Report for Archive"
#
# The localization run moved on from hand-off and is now mid-translation, so
# every "Status" cell that still reads "Ready for handoff" needs to read
# "In Translation" instead:
#   - Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3
#   - zh-cn / de-de sheets: column C ("Status"), rows 2-3
# The shorter replacement text no longer needs as wide a column, so the
# affected "Status" columns are narrowed to match (~12.5 characters).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text -----------------------------------------------
$overview.Range("E2:F2").Value = "In Translation"
$overview.Range("E3:F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the affected "Status" columns -------------------------------
$overview.Range("E:F").ColumnWidth = 12.5
$zhcn.Range("C:C").ColumnWidth = 12.5
$dede.Range("C:C").ColumnWidth = 12.5
